# Update the "K" column (column G) with corrected strikeout counts.
# The sheet previously stored a "Strike#" style value in column G; the
# data has been regenerated so that column G now holds the true "K"
# (strikeouts) value pulled from the source box-score data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 3
    3  = 1
    4  = 1
    5  = 1
    6  = 1
    7  = 0
    8  = 0
    9  = 3
    10 = 0
    11 = 1
    13 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 0
    18 = 0
    19 = 1
    20 = 1
    21 = 3
    23 = 1
    24 = 2
    25 = 1
    26 = 2
    27 = 1
    28 = 2
    29 = 1
    30 = 0
    31 = 1
    32 = 0
    33 = 5
    34 = 1
    35 = 2
    37 = 1
    38 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
